$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (Beta) values for columns F..N
$ws.Range("F2").Value = 302.2197284771133
$ws.Range("G2").Value = 13.90952120653222
$ws.Range("H2").Value = 568.779759737437
$ws.Range("I2").Value = 0.8981474680724735
$ws.Range("J2").Value = 0.01976595230542073
$ws.Range("K2").Value = 1.91562581830158
$ws.Range("L2").Value = 0.2206072381506798
$ws.Range("M2").Value = 0.007285163703287202
$ws.Range("N2").Value = 0.4489508945254752

# Update existing row 3 (Gamma) values for columns F..N
$ws.Range("F3").Value = 0.01004502079645653
$ws.Range("G3").Value = 0.003021358651792475
$ws.Range("H3").Value = 0.01671115311239314
$ws.Range("I3").Value = 0.009354655345478125
$ws.Range("J3").Value = 0.002847080942775661
$ws.Range("K3").Value = 0.01553659064173481
$ws.Range("L3").Value = 0.01008262332407133
$ws.Range("M3").Value = 0.003068835559993827
$ws.Range("N3").Value = 0.01674621405088857

# Add new row 4 (Beta + Gamma) - copy formatting from row 2's A cell first
$ws.Range("A2").Copy($ws.Range("A4"))

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Beta + Gamma"
$ws.Range("C4").Value = 12.00687180793019
$ws.Range("D4").Value = 1.974477778970852
$ws.Range("E4").Value = 0.1537386519519979
$ws.Range("F4").Value = 302.2297734979098
$ws.Range("G4").Value = 13.91254256518401
$ws.Range("H4").Value = 568.7964708905492
$ws.Range("I4").Value = 0.9075021234179516
$ws.Range("J4").Value = 0.02261303324819639
$ws.Range("K4").Value = 1.931162408943315
$ws.Range("L4").Value = 0.2306898614747512
$ws.Range("M4").Value = 0.01035399926328103
$ws.Range("N4").Value = 0.4656971085763638
